$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
# Data through 2021-09-09 (commit: "aggiornamento a 9/09 compreso")
$data = @(
    @(44441, 0, 13, 53.89271204709394),
    @(44442, 9, 19, 78.76627145344499),
    @(44443, 3, 22, 91.20305115662052),
    @(44444, 1, 18, 74.62067821905315),
    @(44445, 6, 23, 95.34864439101236),
    @(44446, 5, 24, 99.4942376254042),
    @(44447, 1, 25, 103.639830859796),
    @(44448, 1, 26, 107.7854240941879)
)

$lastRow = 366
$startRow = $lastRow + 1

# Copy the formatting (style / number format) of the last existing date cell
# so the newly appended date cells in column A keep the same look.
$ws.Cells.Item($lastRow, 1).Copy()

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = $false
